# Add a new "2021" column (column R) to the 13.1.1 disasters-deaths table,
# mirroring the formatting of the existing "2020" column (Q) and filling in
# the new year's data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 2021 values for rows 4-34 (row 3 is a blank/bordered spacer row, so it
# only needs formatting, no value). A value of "-" mirrors the existing
# shared-string placeholder used elsewhere in the sheet for "no data".
$values = @{
    4  = 2021
    5  = 109
    6  = 74
    7  = 35
    8  = 36
    9  = 35
    10 = 1
    11 = 15
    12 = 8
    13 = 7
    14 = 12
    15 = 7
    16 = 5
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 17
    21 = 8
    22 = 9
    23 = 9
    24 = 7
    25 = 2
    26 = 20
    27 = 9
    28 = 11
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

for ($row = 3; $row -le 34; $row++) {
    $src = $ws.Range("Q$row")
    $dst = $ws.Range("R$row")

    # Copy column Q's formatting (border/font/alignment) into column R so the
    # new column visually matches the rest of the table.
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null

    if ($values.ContainsKey($row)) {
        $dst.Value = $values[$row]
    }
}

$excel.CutCopyMode = 0

# Match the author's final selection/active cell.
$ws.Range("R35").Select() | Out-Null
